# daily auto push: 2025-10-05 07:23 UTC
# Append the new daily tracking row (row 64) at the bottom of the sheet's
# data table, mirroring the existing rows (date, weekday, hour, rank count).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 64

# Column A holds the date as literal text (e.g. "2025/10/05"), just like
# every other row already in the sheet. A leading apostrophe forces Excel
# to store it as text instead of auto-converting it to a date serial
# number; Style is then reset to "Normal" so no extra per-cell formatting
# is left behind (matching the unstyled cells used by the other data rows).
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Value = "'2025/10/05"
$dateCell.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "日"
$ws.Cells.Item($newRow, 3).Value = 16
$ws.Cells.Item($newRow, 4).Value = 5
